$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (B8)
$ws.Range("B8").Value = "2024-06-10T07:36:07+00:00"

# Update the Context values (B20, B21)
$ws.Range("B20").Value = "element:MedicationRequest"
$ws.Range("B21").Value = "element:MedicationAdministration"
